# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45924
$ws.Range("B2").Value = 53.51
$ws.Range("C2").Value = 39.8
$ws.Range("D2").Value = 48.6
$ws.Range("E2").Value = 42
$ws.Range("F2").Value = 41.5
$ws.Range("G2").Value = 48.6
$ws.Range("H2").Value = 48.8
$ws.Range("I2").Value = 66.95999999999999
$ws.Range("J2").Value = 69.01000000000001
$ws.Range("K2").Value = 57.91
$ws.Range("L2").Value = 9.52
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 0.65
$ws.Range("O2").Value = 1.55
$ws.Range("P2").Value = 4.31
$ws.Range("Q2").Value = 20.1
$ws.Range("R2").Value = 5.79
$ws.Range("S2").Value = 20.1
$ws.Range("T2").Value = 72.72
$ws.Range("U2").Value = 95
$ws.Range("V2").Value = 116.12
$ws.Range("W2").Value = 105.77
$ws.Range("X2").Value = 90
$ws.Range("Y2").Value = 83.08
$ws.Range("Z2").Value = 47.6
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 98.73999999999999
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 110.94
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 86.54000000000001
$ws.Range("AG2").Value = "1h-17h"
